# Weekly update: insert a new Cilantro price record as the latest row (row 47),
# pushing the existing historical rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 47 (existing rows 47-102 shift down to 48-103).
# Insert copies formatting from the row above, matching the date-format style
# already used by every other row in column D.
$ws.Rows.Item(47).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(47, 1).Value = 8
$ws.Cells.Item(47, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(47, 3).Value = "Coquimbo"
$ws.Cells.Item(47, 4).Value = "12/9/2021"
$ws.Cells.Item(47, 5).Value = 4
$ws.Cells.Item(47, 6).Value = 100112040
$ws.Cells.Item(47, 7).Value = "Cilantro"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 3200
$ws.Cells.Item(47, 11).Value = 1500
$ws.Cells.Item(47, 12).Value = 2000
$ws.Cells.Item(47, 13).Value = 1750
$ws.Cells.Item(47, 14).Value = "`$/atado 1 a 1,5 kilos"
$ws.Cells.Item(47, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(47, 16).Value = 1167
$ws.Cells.Item(47, 17).Value = 1.5
$ws.Cells.Item(47, 18).Value = "Hortaliza"
